$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Actualiza base de datos EC: desplaza los periodos de mora existentes y
# agrega el nuevo periodo (parte 1 de nuevos estado de cuenta).
# Valores anteriores (E16:E19): 2507, 2506, 2505, 2504
# Valores nuevos      (E16:E19): 2505, 2506, 2507, 2508
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2507"
$ws.Range("E19").Value = "2508"
